$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E7").Value = "Lol"
$ws.Range("E9").Value = "rue de la loge 38"
$ws.Range("E10").Value = "7866          Bois de Lessines"
$ws.Range("E13").Value = "Bois de Lessines, 20/11/2023"
$ws.Range("E15").Value = "Facture n°23-002"
